$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.400.55'
$ws.Range('E2').Value = '  -3.98%  '
$ws.Range('D3').Value = '1.570.27'
$ws.Range('E3').Value = '  -3.53%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.001'
$ws.Range('E5').Value = '  +0.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '290.90'
$ws.Range('E6').Value = '  -2.36%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3689'
$ws.Range('E7').Value = '  -2.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '49.32'
$ws.Range('E8').Value = '  -1.36%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3395'
$ws.Range('E9').Value = '  -2.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.170'
$ws.Range('E10').Value = '  -2.60%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07657'
$ws.Range('E11').Value = '  -4.72%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.002'
$ws.Range('E12').Value = '  +0.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.31'
$ws.Range('E13').Value = '  -2.50%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.061'
$ws.Range('E14').Value = '  -3.41%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.928'
$ws.Range('E15').Value = '  -3.75%  '
$ws.Range('D16').Value = '1.572.58'
$ws.Range('E16').Value = '  -3.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001129'
$ws.Range('E17').Value = '  -5.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '90.13'
$ws.Range('E18').Value = '  -4.67%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06750'
$ws.Range('E19').Value = '  -2.75%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  -0.04%  '
$ws.Range('E21').Value = '  -4.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '16.55'
$ws.Range('E22').Value = '  -3.82%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.5310'
$ws.Range('E23').Value = '  -7.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.02'
$ws.Range('E24').Value = '  -2.62%  '
$ws.Range('D25').Value = '22.406.40'
$ws.Range('E25').Value = '  -4.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.362'
$ws.Range('E26').Value = '  -2.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.864'
$ws.Range('E27').Value = '  -2.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.10'
$ws.Range('E28').Value = '  -3.37%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '146.37'
$ws.Range('E29').Value = '  -1.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.980'
$ws.Range('E30').Value = '  -3.49%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '125.80'
$ws.Range('E31').Value = '  -3.68%  '
$ws.Range('D32').Value = '1.744.60'
$ws.Range('E32').Value = '  -3.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.021'
$ws.Range('E33').Value = '  +4.00%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.206'
$ws.Range('E34').Value = '  -7.53%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.04'
$ws.Range('E36').Value = '  -9.61%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.08478'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02545'
$ws.Range('E38').Value = '  -3.89%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2322'
$ws.Range('E39').Value = '  -3.43%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.524'
$ws.Range('E40').Value = '  -4.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.06464'
$ws.Range('E41').Value = '  -4.28%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.297'
$ws.Range('E42').Value = '  +0.89%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.70'
$ws.Range('E43').Value = '  -7.69%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6349'
$ws.Range('E44').Value = '  -6.25%  '
$ws.Range('E45').Value = '  -7.33%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.000'
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5988'
$ws.Range('E47').Value = '  -4.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.764'
$ws.Range('E48').Value = '  -3.16%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.113'
$ws.Range('E49').Value = '  -4.85%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.265'
$ws.Range('E50').Value = '  +3.73%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '124.81'
$ws.Range('E51').Value = '  -0.98%  '
